$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 27780020  # ALC!H40
$ws.Cells.Item(40, 9).Value = 1950  # ALC!I40
$ws.Cells.Item(40, 11).Value = 1950  # ALC!K40
$ws.Cells.Item(40, 13).Value = -1775  # ALC!M40

$ws.Cells.Item(116, 8).Value = 4459.6  # ALC!H116
$ws.Cells.Item(116, 9).Value = 4881.6665  # ALC!I116
$ws.Cells.Item(116, 11).Value = 4881.6665  # ALC!K116
$ws.Cells.Item(116, 13).Value = -1439.6665  # ALC!M116

$ws.Cells.Item(118, 8).Value = 368.33334  # ALC!H118
$ws.Cells.Item(118, 9).Value = 242  # ALC!I118
$ws.Cells.Item(118, 11).Value = 726  # ALC!K118
$ws.Cells.Item(118, 13).Value = 931  # ALC!M118

$ws.Cells.Item(135, 8).Value = 2357.6316  # ALC!H135
$ws.Cells.Item(135, 9).Value = 1342.5714  # ALC!I135
$ws.Cells.Item(135, 11).Value = 12083.1426  # ALC!K135
$ws.Cells.Item(135, 13).Value = -9548.142600000001  # ALC!M135

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2730.0557  # ARM!H2
$ws.Cells.Item(2, 9).Value = 2574.3572  # ARM!I2
$ws.Cells.Item(2, 10).Value = 3275  # ARM!J2
$ws.Cells.Item(2, 11).Value = 2574.3572  # ARM!K2
$ws.Cells.Item(2, 12).Value = 3275  # ARM!L2
$ws.Cells.Item(2, 13).Value = -2461.3572  # ARM!M2
$ws.Cells.Item(2, 14).Value = -3501  # ARM!N2

$ws.Cells.Item(4, 8).Value = 269  # ARM!H4
$ws.Cells.Item(4, 9).Value = 261.5  # ARM!I4
$ws.Cells.Item(4, 10).Value = 278  # ARM!J4
$ws.Cells.Item(4, 11).Value = 261.5  # ARM!K4
$ws.Cells.Item(4, 12).Value = 278  # ARM!L4
$ws.Cells.Item(4, 13).Value = -145.5  # ARM!M4
$ws.Cells.Item(4, 14).Value = -510  # ARM!N4

$ws.Cells.Item(5, 8).Value = 60  # ARM!H5
$ws.Cells.Item(5, 9).Value = 24.5  # ARM!I5
$ws.Cells.Item(5, 10).Value = 83.666664  # ARM!J5
$ws.Cells.Item(5, 11).Value = 24.5  # ARM!K5
$ws.Cells.Item(5, 12).Value = 83.666664  # ARM!L5
$ws.Cells.Item(5, 13).Value = 87.5  # ARM!M5
$ws.Cells.Item(5, 14).Value = -307.666664  # ARM!N5

$ws.Cells.Item(32, 8).Value = 2610590  # ARM!H32
$ws.Cells.Item(32, 9).Value = 5424.171  # ARM!I32
$ws.Cells.Item(32, 10).Value = 17869418  # ARM!J32
$ws.Cells.Item(32, 11).Value = 5424.171  # ARM!K32
$ws.Cells.Item(32, 12).Value = 17869418  # ARM!L32
$ws.Cells.Item(32, 13).Value = -5137.171  # ARM!M32
$ws.Cells.Item(32, 14).Value = -17869992  # ARM!N32

$ws.Cells.Item(45, 8).Value = 1529.3684  # ARM!H45
$ws.Cells.Item(45, 9).Value = 1472.6154  # ARM!I45
$ws.Cells.Item(45, 10).Value = 1652.3334  # ARM!J45
$ws.Cells.Item(45, 11).Value = 1472.6154  # ARM!K45
$ws.Cells.Item(45, 12).Value = 1652.3334  # ARM!L45
$ws.Cells.Item(45, 13).Value = -1095.6154  # ARM!M45
$ws.Cells.Item(45, 14).Value = -2406.3334  # ARM!N45

$ws.Cells.Item(74, 8).Value = 978.119  # ARM!H74
$ws.Cells.Item(74, 9).Value = 960.71875  # ARM!I74
$ws.Cells.Item(74, 10).Value = 1033.8  # ARM!J74
$ws.Cells.Item(74, 11).Value = 960.71875  # ARM!K74
$ws.Cells.Item(74, 12).Value = 1033.8  # ARM!L74
$ws.Cells.Item(74, 13).Value = -86.71875  # ARM!M74
$ws.Cells.Item(74, 14).Value = -2781.8  # ARM!N74

$ws.Cells.Item(77, 8).Value = 978.119  # ARM!H77
$ws.Cells.Item(77, 9).Value = 960.71875  # ARM!I77
$ws.Cells.Item(77, 10).Value = 1033.8  # ARM!J77
$ws.Cells.Item(77, 11).Value = 4803.59375  # ARM!K77
$ws.Cells.Item(77, 12).Value = 5169  # ARM!L77
$ws.Cells.Item(77, 13).Value = -435.59375  # ARM!M77
$ws.Cells.Item(77, 14).Value = -13905  # ARM!N77

$ws.Cells.Item(102, 8).Value = 3652.4167  # ARM!H102
$ws.Cells.Item(102, 9).Value = 3203.2222  # ARM!I102
$ws.Cells.Item(102, 10).Value = 5000  # ARM!J102
$ws.Cells.Item(102, 11).Value = 3203.2222  # ARM!K102
$ws.Cells.Item(102, 12).Value = 5000  # ARM!L102
$ws.Cells.Item(102, 13).Value = -1581.2222  # ARM!M102
$ws.Cells.Item(102, 14).Value = -8244  # ARM!N102

$ws.Cells.Item(116, 8).Value = 2730.0557  # ARM!H116
$ws.Cells.Item(116, 9).Value = 2574.3572  # ARM!I116
$ws.Cells.Item(116, 10).Value = 3275  # ARM!J116
$ws.Cells.Item(116, 11).Value = 2574.3572  # ARM!K116
$ws.Cells.Item(116, 12).Value = 3275  # ARM!L116
$ws.Cells.Item(116, 13).Value = -280.3571999999999  # ARM!M116
$ws.Cells.Item(116, 14).Value = -7863  # ARM!N116

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2730.0557  # BSM!H3
$ws.Cells.Item(3, 9).Value = 2574.3572  # BSM!I3
$ws.Cells.Item(3, 10).Value = 3275  # BSM!J3
$ws.Cells.Item(3, 11).Value = 2574.3572  # BSM!K3
$ws.Cells.Item(3, 12).Value = 3275  # BSM!L3
$ws.Cells.Item(3, 13).Value = -2460.3572  # BSM!M3
$ws.Cells.Item(3, 14).Value = -3503  # BSM!N3

$ws.Cells.Item(4, 8).Value = 60  # BSM!H4
$ws.Cells.Item(4, 9).Value = 24.5  # BSM!I4
$ws.Cells.Item(4, 10).Value = 83.666664  # BSM!J4
$ws.Cells.Item(4, 11).Value = 24.5  # BSM!K4
$ws.Cells.Item(4, 12).Value = 83.666664  # BSM!L4
$ws.Cells.Item(4, 13).Value = 90.5  # BSM!M4
$ws.Cells.Item(4, 14).Value = -313.666664  # BSM!N4

$ws.Cells.Item(80, 8).Value = 536.625  # BSM!H80
$ws.Cells.Item(80, 9).Value = 281  # BSM!I80
$ws.Cells.Item(80, 10).Value = 690  # BSM!J80
$ws.Cells.Item(80, 11).Value = 281  # BSM!K80
$ws.Cells.Item(80, 12).Value = 690  # BSM!L80
$ws.Cells.Item(80, 13).Value = 717  # BSM!M80
$ws.Cells.Item(80, 14).Value = -2686  # BSM!N80

$ws.Cells.Item(83, 8).Value = 536.625  # BSM!H83
$ws.Cells.Item(83, 9).Value = 281  # BSM!I83
$ws.Cells.Item(83, 10).Value = 690  # BSM!J83
$ws.Cells.Item(83, 11).Value = 1405  # BSM!K83
$ws.Cells.Item(83, 12).Value = 3450  # BSM!L83
$ws.Cells.Item(83, 13).Value = 3587  # BSM!M83
$ws.Cells.Item(83, 14).Value = -13434  # BSM!N83

$ws.Cells.Item(86, 8).Value = 3062.1875  # BSM!H86
$ws.Cells.Item(86, 9).Value = 2199.5454  # BSM!I86
$ws.Cells.Item(86, 10).Value = 4960  # BSM!J86
$ws.Cells.Item(86, 11).Value = 2199.5454  # BSM!K86
$ws.Cells.Item(86, 12).Value = 4960  # BSM!L86
$ws.Cells.Item(86, 13).Value = -1076.5454  # BSM!M86
$ws.Cells.Item(86, 14).Value = -7206  # BSM!N86

$ws.Cells.Item(89, 8).Value = 3062.1875  # BSM!H89
$ws.Cells.Item(89, 9).Value = 2199.5454  # BSM!I89
$ws.Cells.Item(89, 10).Value = 4960  # BSM!J89
$ws.Cells.Item(89, 11).Value = 10997.727  # BSM!K89
$ws.Cells.Item(89, 12).Value = 24800  # BSM!L89
$ws.Cells.Item(89, 13).Value = -5381.726999999999  # BSM!M89
$ws.Cells.Item(89, 14).Value = -36032  # BSM!N89

$ws.Cells.Item(94, 8).Value = 935.4  # BSM!H94
$ws.Cells.Item(94, 9).Value = 625.46155  # BSM!I94
$ws.Cells.Item(94, 11).Value = 625.46155  # BSM!K94
$ws.Cells.Item(94, 13).Value = -174.46155  # BSM!M94

$ws.Cells.Item(99, 8).Value = 45456816  # BSM!H99
$ws.Cells.Item(99, 9).Value = 50002256  # BSM!I99
$ws.Cells.Item(99, 10).Value = 2411  # BSM!J99
$ws.Cells.Item(99, 11).Value = 50002256  # BSM!K99
$ws.Cells.Item(99, 12).Value = 2411  # BSM!L99
$ws.Cells.Item(99, 13).Value = -50000758  # BSM!M99
$ws.Cells.Item(99, 14).Value = -5407  # BSM!N99

$ws.Cells.Item(105, 8).Value = 5312.4116  # BSM!H105
$ws.Cells.Item(105, 9).Value = 4608.3335  # BSM!I105
$ws.Cells.Item(105, 10).Value = 7002.2  # BSM!J105
$ws.Cells.Item(105, 11).Value = 4608.3335  # BSM!K105
$ws.Cells.Item(105, 12).Value = 7002.2  # BSM!L105
$ws.Cells.Item(105, 13).Value = -2861.3335  # BSM!M105
$ws.Cells.Item(105, 14).Value = -10496.2  # BSM!N105

$ws.Cells.Item(134, 8).Value = 29246.41  # BSM!H134
$ws.Cells.Item(134, 9).Value = 4531.722  # BSM!I134
$ws.Cells.Item(134, 10).Value = 50430.43  # BSM!J134
$ws.Cells.Item(134, 11).Value = 13595.166  # BSM!K134
$ws.Cells.Item(134, 12).Value = 151291.29  # BSM!L134
$ws.Cells.Item(134, 13).Value = -11060.166  # BSM!M134
$ws.Cells.Item(134, 14).Value = -156361.29  # BSM!N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 151.81818  # CRP!H7
$ws.Cells.Item(7, 9).Value = 147  # CRP!I7
$ws.Cells.Item(7, 11).Value = 147  # CRP!K7
$ws.Cells.Item(7, 13).Value = -34  # CRP!M7

$ws.Cells.Item(31, 8).Value = 2733.2083  # CRP!H31
$ws.Cells.Item(31, 10).Value = 1457.2  # CRP!J31
$ws.Cells.Item(31, 12).Value = 1457.2  # CRP!L31
$ws.Cells.Item(31, 14).Value = -2047.2  # CRP!N31

$ws.Cells.Item(34, 8).Value = 2733.2083  # CRP!H34
$ws.Cells.Item(34, 10).Value = 1457.2  # CRP!J34
$ws.Cells.Item(34, 12).Value = 1457.2  # CRP!L34
$ws.Cells.Item(34, 14).Value = -1861.2  # CRP!N34

$ws.Cells.Item(105, 8).Value = 1470.625  # CRP!H105
$ws.Cells.Item(105, 9).Value = 1395  # CRP!I105
$ws.Cells.Item(105, 11).Value = 1395  # CRP!K105
$ws.Cells.Item(105, 13).Value = 352  # CRP!M105

$ws.Cells.Item(132, 8).Value = 2653.4285  # CRP!H132
$ws.Cells.Item(132, 9).Value = 2135.2917  # CRP!I132
$ws.Cells.Item(132, 10).Value = 3783.9092  # CRP!J132
$ws.Cells.Item(132, 11).Value = 6405.875100000001  # CRP!K132
$ws.Cells.Item(132, 12).Value = 11351.7276  # CRP!L132
$ws.Cells.Item(132, 13).Value = -3875.875100000001  # CRP!M132
$ws.Cells.Item(132, 14).Value = -16411.7276  # CRP!N132

$ws.Cells.Item(134, 8).Value = 2750.6924  # CRP!H134
$ws.Cells.Item(134, 9).Value = 2091.3333  # CRP!I134
$ws.Cells.Item(134, 11).Value = 6273.999899999999  # CRP!K134
$ws.Cells.Item(134, 13).Value = -3738.999899999999  # CRP!M134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 11327.272  # CUL!H87
$ws.Cells.Item(87, 9).Value = 2916.6667  # CUL!I87
$ws.Cells.Item(87, 10).Value = 21420  # CUL!J87
$ws.Cells.Item(87, 11).Value = 8750.000100000001  # CUL!K87
$ws.Cells.Item(87, 12).Value = 64260  # CUL!L87
$ws.Cells.Item(87, 13).Value = -7502.000100000001  # CUL!M87
$ws.Cells.Item(87, 14).Value = -66756  # CUL!N87

$ws.Cells.Item(90, 8).Value = 11327.272  # CUL!H90
$ws.Cells.Item(90, 9).Value = 2916.6667  # CUL!I90
$ws.Cells.Item(90, 10).Value = 21420  # CUL!J90
$ws.Cells.Item(90, 11).Value = 26250.0003  # CUL!K90
$ws.Cells.Item(90, 12).Value = 192780  # CUL!L90
$ws.Cells.Item(90, 13).Value = -20010.0003  # CUL!M90
$ws.Cells.Item(90, 14).Value = -205260  # CUL!N90

$ws.Cells.Item(137, 8).Value = 45691.375  # CUL!H137
$ws.Cells.Item(137, 9).Value = 1984  # CUL!I137
$ws.Cells.Item(137, 11).Value = 5952  # CUL!K137
$ws.Cells.Item(137, 13).Value = -852  # CUL!M137

$ws.Cells.Item(140, 8).Value = 4438.9  # CUL!H140
$ws.Cells.Item(140, 9).Value = 4662.5  # CUL!I140
$ws.Cells.Item(140, 10).Value = 3544.5  # CUL!J140
$ws.Cells.Item(140, 11).Value = 13987.5  # CUL!K140
$ws.Cells.Item(140, 12).Value = 10633.5  # CUL!L140
$ws.Cells.Item(140, 13).Value = -8807.5  # CUL!M140
$ws.Cells.Item(140, 14).Value = -20993.5  # CUL!N140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 61.07143  # GSM!H2
$ws.Cells.Item(2, 9).Value = 33.75  # GSM!I2
$ws.Cells.Item(2, 11).Value = 33.75  # GSM!K2
$ws.Cells.Item(2, 13).Value = 79.25  # GSM!M2

$ws.Cells.Item(97, 8).Value = 1985.7142  # GSM!H97
$ws.Cells.Item(97, 9).Value = 1985.7142  # GSM!I97
$ws.Cells.Item(97, 10).Value = 0  # GSM!J97
$ws.Cells.Item(97, 11).Value = 1985.7142  # GSM!K97
$ws.Cells.Item(97, 12).Value = 0  # GSM!L97
$ws.Cells.Item(97, 13).Value = -1489.7142  # GSM!M97
$ws.Cells.Item(97, 14).Value = $null  # GSM!N97

$ws.Cells.Item(102, 8).Value = 1121.4517  # GSM!H102
$ws.Cells.Item(102, 9).Value = 867.88464  # GSM!I102
$ws.Cells.Item(102, 10).Value = 2440  # GSM!J102
$ws.Cells.Item(102, 11).Value = 867.88464  # GSM!K102
$ws.Cells.Item(102, 12).Value = 2440  # GSM!L102
$ws.Cells.Item(102, 13).Value = 754.11536  # GSM!M102
$ws.Cells.Item(102, 14).Value = -5684  # GSM!N102

$ws.Cells.Item(132, 8).Value = 2862.2  # GSM!H132
$ws.Cells.Item(132, 9).Value = 941.3333  # GSM!I132
$ws.Cells.Item(132, 10).Value = 3685.4285  # GSM!J132
$ws.Cells.Item(132, 11).Value = 2823.9999  # GSM!K132
$ws.Cells.Item(132, 12).Value = 11056.2855  # GSM!L132
$ws.Cells.Item(132, 13).Value = -293.9998999999998  # GSM!M132
$ws.Cells.Item(132, 14).Value = -16116.2855  # GSM!N132

$ws.Cells.Item(135, 8).Value = 40126.668  # GSM!H135
$ws.Cells.Item(135, 10).Value = 40126.668  # GSM!J135
$ws.Cells.Item(135, 12).Value = 40126.668  # GSM!L135
$ws.Cells.Item(135, 14).Value = -50266.668  # GSM!N135

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1446.95  # LTW!H22
$ws.Cells.Item(22, 9).Value = 870  # LTW!I22
$ws.Cells.Item(22, 10).Value = 1639.2667  # LTW!J22
$ws.Cells.Item(22, 11).Value = 870  # LTW!K22
$ws.Cells.Item(22, 12).Value = 1639.2667  # LTW!L22
$ws.Cells.Item(22, 13).Value = -575  # LTW!M22
$ws.Cells.Item(22, 14).Value = -2229.2667  # LTW!N22

$ws.Cells.Item(27, 8).Value = 1446.95  # LTW!H27
$ws.Cells.Item(27, 9).Value = 870  # LTW!I27
$ws.Cells.Item(27, 10).Value = 1639.2667  # LTW!J27
$ws.Cells.Item(27, 11).Value = 870  # LTW!K27
$ws.Cells.Item(27, 12).Value = 1639.2667  # LTW!L27
$ws.Cells.Item(27, 13).Value = -763  # LTW!M27
$ws.Cells.Item(27, 14).Value = -1853.2667  # LTW!N27

$ws.Cells.Item(55, 8).Value = 438.42856  # LTW!H55
$ws.Cells.Item(55, 9).Value = 440  # LTW!I55
$ws.Cells.Item(55, 10).Value = 437.25  # LTW!J55
$ws.Cells.Item(55, 11).Value = 440  # LTW!K55
$ws.Cells.Item(55, 12).Value = 437.25  # LTW!L55
$ws.Cells.Item(55, 13).Value = -267  # LTW!M55
$ws.Cells.Item(55, 14).Value = -783.25  # LTW!N55

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 610.7931  # WVR!H107
$ws.Cells.Item(107, 9).Value = 678.4091  # WVR!I107
$ws.Cells.Item(107, 10).Value = 398.2857  # WVR!J107
$ws.Cells.Item(107, 11).Value = 2035.2273  # WVR!K107
$ws.Cells.Item(107, 12).Value = 1194.8571  # WVR!L107
$ws.Cells.Item(107, 13).Value = -115.2273  # WVR!M107
$ws.Cells.Item(107, 14).Value = -5034.8571  # WVR!N107
